# Update the "Estado de Cuenta" (EC) database with the new statement rows.
# This reflects interleaving/recalculating the periods and values for
# trabajadores DIANA CAROLINA DE AVILA TORDECILLA (CC 1128046581) and
# KATIA DEL PILAR TOVAR TINOCO (CC 45514830), rows 16-73 of Hoja1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$data = @(
    @(16, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "2001", 120000, 3000000),
    @(17, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1912", 120000, 3000000),
    @(18, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1911", 120000, 3000000),
    @(19, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1910", 120000, 3000000),
    @(20, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1909", 120000, 3000000),
    @(21, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1908", 120000, 3000000),
    @(22, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1907", 120000, 3000000),
    @(23, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1906", 120000, 3000000),
    @(24, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1905", 120000, 3000000),
    @(25, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1904", 120000, 3000000),
    @(26, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1903", 120000, 3000000),
    @(27, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1902", 120000, 3000000),
    @(28, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1901", 120000, 3000000),
    @(29, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1812", 120000, 3000000),
    @(30, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1811", 120000, 3000000),
    @(31, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1810", 120000, 3000000),
    @(32, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1809", 120000, 3000000),
    @(33, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1808", 120000, 3000000),
    @(34, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1807", 120000, 3000000),
    @(35, "45514830", "KATIA DEL PILAR TOVAR TINOCO", "2001", 120000, 3000000),
    @(36, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1807", 120000, 3122700),
    @(37, "45514830", "KATIA DEL PILAR TOVAR TINOCO", "1807", 120000, 3000000),
    @(38, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1808", 120000, 3122700),
    @(39, "45514830", "KATIA DEL PILAR TOVAR TINOCO", "1808", 120000, 3000000),
    @(40, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1809", 124908, 3122700),
    @(41, "45514830", "KATIA DEL PILAR TOVAR TINOCO", "1809", 120000, 3000000),
    @(42, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1810", 124908, 3122700),
    @(43, "45514830", "KATIA DEL PILAR TOVAR TINOCO", "1810", 120000, 3000000),
    @(44, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1811", 124908, 3122700),
    @(45, "45514830", "KATIA DEL PILAR TOVAR TINOCO", "1811", 120000, 3000000),
    @(46, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1812", 124908, 3122700),
    @(47, "45514830", "KATIA DEL PILAR TOVAR TINOCO", "1812", 120000, 3000000),
    @(48, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1901", 124908, 3122700),
    @(49, "45514830", "KATIA DEL PILAR TOVAR TINOCO", "1901", 120000, 3000000),
    @(50, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1902", 124908, 3122700),
    @(51, "45514830", "KATIA DEL PILAR TOVAR TINOCO", "1902", 120000, 3000000),
    @(52, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1903", 124908, 3122700),
    @(53, "45514830", "KATIA DEL PILAR TOVAR TINOCO", "1903", 120000, 3000000),
    @(54, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1904", 124908, 3122700),
    @(55, "45514830", "KATIA DEL PILAR TOVAR TINOCO", "1904", 120000, 3000000),
    @(56, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1905", 124908, 3122700),
    @(57, "45514830", "KATIA DEL PILAR TOVAR TINOCO", "1905", 120000, 3000000),
    @(58, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1906", 124908, 3122700),
    @(59, "45514830", "KATIA DEL PILAR TOVAR TINOCO", "1906", 120000, 3000000),
    @(60, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1907", 124908, 3122700),
    @(61, "45514830", "KATIA DEL PILAR TOVAR TINOCO", "1907", 120000, 3000000),
    @(62, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1908", 124908, 3122700),
    @(63, "45514830", "KATIA DEL PILAR TOVAR TINOCO", "1908", 120000, 3000000),
    @(64, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1909", 124908, 3122700),
    @(65, "45514830", "KATIA DEL PILAR TOVAR TINOCO", "1909", 120000, 3000000),
    @(66, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1910", 124908, 3122700),
    @(67, "45514830", "KATIA DEL PILAR TOVAR TINOCO", "1910", 120000, 3000000),
    @(68, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1911", 124908, 3122700),
    @(69, "45514830", "KATIA DEL PILAR TOVAR TINOCO", "1911", 120000, 3000000),
    @(70, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "1912", 124908, 3122700),
    @(71, "45514830", "KATIA DEL PILAR TOVAR TINOCO", "1912", 120000, 3000000),
    @(72, "1128046581", "DIANA CAROLINA DE AVILA TORDECILLA", "2001", 87436, 3122700),
    @(73, "45514830", "KATIA DEL PILAR TOVAR TINOCO", "1611", 84000, 3000000),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}

